$wb = $excel.ActiveWorkbook

$wsTotalEnergy = $wb.Worksheets.Item("TotalEnergy")
$wsMeters      = $wb.Worksheets.Item("Meters")
$wsSources     = $wb.Worksheets.Item("Sources")

# --- TotalEnergy sheet: the extra output rows (Natural Gas / Electricity
#     totals, Electricity Heat Rejection) are no longer needed - only the
#     header and "Total Site Energy" row remain.
$wsTotalEnergy.Rows("3:6").Delete()
$wsTotalEnergy.Range("A3:XFD18").Select()

# --- Meters sheet: fixed timestep bug - "timestep" was being looked up as
#     a report frequency ("Monthly") instead of the literal "Timestep"
#     meter-reporting keyword, and the extra Gas:Facility meter row is
#     removed as part of the uncertainty-code cleanup.
$wsMeters.Range("B2").Value = "Timestep"
$wsMeters.Rows("3:3").Delete()
$wsMeters.Range("B2").Select()

# --- Sources sheet: selection moved.
$wsSources.Range("E8").Select()

# --- Meters is the sheet left active/visible when the workbook was saved.
$wsMeters.Activate()
